$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 1.94
$ws.Range("R2").Value = 1.79
$ws.Range("R3").Value = 1.5
$ws.Range("V3").Value = 1.54
$ws.Range("R4").Value = 1.33
$ws.Range("V4").Value = 1.47
$ws.Range("S8").Value = 1.37
$ws.Range("G12").Value = 1.75
$ws.Range("I12").Value = 4.15
$ws.Range("J12").Value = 2.27
$ws.Range("K12").Value = 2.2
$ws.Range("L12").Value = 4.45
$ws.Range("W12").Value = 8.25
$ws.Range("X12").Value = 9.25
$ws.Range("AA12").Value = 13
$ws.Range("AD12").Value = 7.2
$ws.Range("AE12").Value = 14
$ws.Range("AG12").Value = 13
$ws.Range("AH12").Value = 25
$ws.Range("AI12").Value = 13.5
$ws.Range("AJ12").Value = 65
$ws.Range("AK12").Value = 37
$ws.Range("AM12").Value = 400
$ws.Range("AN12").Value = 3.7
$ws.Range("AO12").Value = 8.5
$ws.Range("AT12").Value = 2.9
$ws.Range("AU12").Value = 7
$ws.Range("AW12").Value = 6
$ws.Range("AX12").Value = 23
$ws.Range("AY12").Value = 27
$ws.Range("AZ12").Value = 120
$ws.Range("S16").Value = 1.41
$ws.Range("T16").Value = 2.62
$ws.Range("S17").Value = 1.41
$ws.Range("T17").Value = 2.62
$ws.Range("S21").Value = 1.41
$ws.Range("T21").Value = 2.62
$ws.Range("S25").Value = 1.41
$ws.Range("T25").Value = 2.62
$ws.Range("S26").Value = 1.19
$ws.Range("G29").Value = 2.5
$ws.Range("H29").Value = 3.55
$ws.Range("K29").Value = 2.22
$ws.Range("M29").Value = 1.04
$ws.Range("N29").Value = 8.25
$ws.Range("O29").Value = 1.24
$ws.Range("P29").Value = 3.65
$ws.Range("Q29").Value = 1.75
$ws.Range("R29").Value = 2.02
$ws.Range("S29").Value = 1.35
$ws.Range("T29").Value = 2.92
$ws.Range("U29").Value = 1.62
$ws.Range("V29").Value = 2.15
$ws.Range("W29").Value = 9.75
$ws.Range("X29").Value = 13.5
$ws.Range("AA29").Value = 19.5
$ws.Range("AB29").Value = 26
$ws.Range("AC29").Value = 8.25
$ws.Range("AD29").Value = 7
$ws.Range("AE29").Value = 13
$ws.Range("AF29").Value = 50
$ws.Range("AG29").Value = 9.5
$ws.Range("AH29").Value = 13
$ws.Range("AK29").Value = 18.5
$ws.Range("AL29").Value = 26
$ws.Range("AM29").Value = 350
$ws.Range("AN29").Value = 4.6
$ws.Range("AO29").Value = 13
$ws.Range("AP29").Value = 20
$ws.Range("AR29").Value = 80
$ws.Range("AT29").Value = 2.92
$ws.Range("AU29").Value = 7
$ws.Range("AV29").Value = 55
$ws.Range("AW29").Value = 4.5
$ws.Range("AY29").Value = 19.5
$ws.Range("J35").Value = 2.32
$ws.Range("N35").Value = 8
$ws.Range("O35").Value = 1.32
$ws.Range("P35").Value = 2.87
$ws.Range("Q35").Value = 1.93
$ws.Range("R35").Value = 1.7
$ws.Range("U35").Value = 1.82
$ws.Range("V35").Value = 1.78
$ws.Range("W35").Value = 6.5
$ws.Range("X35").Value = 8
$ws.Range("AA35").Value = 15
$ws.Range("AB35").Value = 29
$ws.Range("AC35").Value = 9
$ws.Range("AD35").Value = 6.6
$ws.Range("AG35").Value = 11
$ws.Range("AP35").Value = 17.5
$ws.Range("AR35").Value = 60
$ws.Range("AT35").Value = 2.6
$ws.Range("AZ35").Value = 120
$ws.Range("G38").Value = 10.25
$ws.Range("O38").Value = 1.23
$ws.Range("P38").Value = 3.75
$ws.Range("S38").Value = 1.34
$ws.Range("T38").Value = 3
$ws.Range("U38").Value = 2.35
$ws.Range("V38").Value = 1.53
$ws.Range("W38").Value = 24
$ws.Range("Z38").Value = 350
$ws.Range("AI38").Value = 9.5
$ws.Range("AT38").Value = 3
$ws.Range("AU38").Value = 10.25
$ws.Range("AV38").Value = 120
$ws.Range("AY38").Value = 18.5
$ws.Range("BA38").Value = 50
$ws.Range("G39").Value = 1.4
$ws.Range("H39").Value = 4.4
$ws.Range("I39").Value = 6.2
$ws.Range("J39").Value = 1.87
$ws.Range("N39").Value = 9.25
$ws.Range("O39").Value = 1.17
$ws.Range("P39").Value = 4.45
$ws.Range("Q39").Value = 1.52
$ws.Range("R39").Value = 2.37
$ws.Range("U39").Value = 1.7
$ws.Range("V39").Value = 2.02
$ws.Range("W39").Value = 8.75
$ws.Range("X39").Value = 7.9
$ws.Range("Z39").Value = 10.25
$ws.Range("AA39").Value = 10.5
$ws.Range("AB39").Value = 21
$ws.Range("AC39").Value = 9.25
$ws.Range("AE39").Value = 16.5
$ws.Range("AF39").Value = 60
$ws.Range("AK39").Value = 60
$ws.Range("AL39").Value = 50
$ws.Range("AM39").Value = 400
$ws.Range("AN39").Value = 3.45
$ws.Range("AO39").Value = 6.4
$ws.Range("AP39").Value = 14
$ws.Range("AR39").Value = 37
$ws.Range("AS39").Value = 150
$ws.Range("AU39").Value = 7.6
$ws.Range("J40").Value = 2.42
$ws.Range("P40").Value = 4.35
$ws.Range("R40").Value = 2.32
$ws.Range("T40").Value = 3.3
$ws.Range("U40").Value = 1.5
$ws.Range("V40").Value = 2.4
$ws.Range("Z40").Value = 18
$ws.Range("AE40").Value = 11.75
$ws.Range("AG40").Value = 14.5
$ws.Range("AH40").Value = 22
$ws.Range("AL40").Value = 26
$ws.Range("AM40").Value = 200
$ws.Range("AO40").Value = 9.5
$ws.Range("AP40").Value = 15
$ws.Range("AR40").Value = 50
$ws.Range("AT40").Value = 3.3
$ws.Range("AX40").Value = 17.5
$ws.Range("AY40").Value = 20
$ws.Range("BA40").Value = 90
$ws.Range("G41").Value = 3.65
$ws.Range("H41").Value = 3.4
$ws.Range("K41").Value = 2.12
$ws.Range("L41").Value = 2.5
$ws.Range("N41").Value = 7.5
$ws.Range("Q41").Value = 1.85
$ws.Range("T41").Value = 2.75
$ws.Range("W41").Value = 11.25
$ws.Range("X41").Value = 21
$ws.Range("AC41").Value = 7.5
$ws.Range("AD41").Value = 6.7
$ws.Range("AG41").Value = 7.4
$ws.Range("AP41").Value = 27
$ws.Range("AT41").Value = 2.75
$ws.Range("AU41").Value = 7.2
$ws.Range("AX41").Value = 9.75
$ws.Range("AY41").Value = 18.5
$ws.Range("AZ41").Value = 35
$ws.Range("G42").Value = 1.88
$ws.Range("H42").Value = 3.45
$ws.Range("I42").Value = 3.55
$ws.Range("J42").Value = 2.45
$ws.Range("K42").Value = 2.18
$ws.Range("L42").Value = 4.05
$ws.Range("O42").Value = 1.32
$ws.Range("P42").Value = 3.1
$ws.Range("Q42").Value = 1.98
$ws.Range("R42").Value = 1.78
$ws.Range("S42").Value = 1.38
$ws.Range("T42").Value = 2.8
$ws.Range("U42").Value = 1.85
$ws.Range("V42").Value = 1.85
$ws.Range("W42").Value = 6.9
$ws.Range("X42").Value = 8.75
$ws.Range("Z42").Value = 15.5
$ws.Range("AA42").Value = 15.5
$ws.Range("AB42").Value = 29
$ws.Range("AD42").Value = 6.8
$ws.Range("AE42").Value = 16
$ws.Range("AT42").Value = 2.8
$ws.Range("AU42").Value = 7.4
$ws.Range("AV42").Value = 70
$ws.Range("AW42").Value = 5.5
$ws.Range("AX42").Value = 20
$ws.Range("AY42").Value = 27
$ws.Range("AZ42").Value = 100
$ws.Range("BA42").Value = 150
$ws.Range("AF42").Value = 80
$ws.Range("AG42").Value = 10.25
$ws.Range("AH42").Value = 19
$ws.Range("AI42").Value = 12.5
$ws.Range("AJ42").Value = 50
$ws.Range("AK42").Value = 35
$ws.Range("AL42").Value = 45
$ws.Range("AM42").Value = 700
$ws.Range("AN42").Value = 3.75
$ws.Range("AO42").Value = 9.5
$ws.Range("AP42").Value = 18.5
$ws.Range("AQ42").Value = 35
$ws.Range("AR42").Value = 65
